$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: update quantities/amounts ---
$ws.Range("B5").Value = 100
$ws.Range("C5").Value = 1500
$ws.Range("F5").Value = 150
$ws.Range("G5").Value = 250
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 250
$ws.Range("L5").Value = 2000

# Give row 5 a thin box border around every cell (matches the "data row" frame)
$ws.Range("A5:L5").Borders.LineStyle = 1

# --- Row 6: wipe out the old totals row, leave it blank ---
$ws.Range("A6:L6").ClearContents()

# --- Row 7: repeat of the header row, framed, bold, first cell blank (single space) ---
$ws.Range("A7").Value = " "
$ws.Range("B7").Value = "8532 nos"
$ws.Range("C7").Value = "8532 amt"
$ws.Range("D7").Value = "8533 nos"
$ws.Range("E7").Value = "8533 amt"
$ws.Range("F7").Value = "8536 nos"
$ws.Range("G7").Value = "8535 amt"
$ws.Range("H7").Value = "8541 nos"
$ws.Range("I7").Value = "8541 amt"
$ws.Range("J7").Value = "8542 nos"
$ws.Range("K7").Value = "8542 amt"
$ws.Range("L7").Value = "Total"

$ws.Range("A7:L7").Borders.LineStyle = 1
$ws.Range("A7:L7").Font.Bold = $true

# --- Rows 8-10: TOTAL / CGST / SGST frame labels ---
$ws.Range("A8").Value = "TOTAL"
$ws.Range("A8").Borders.LineStyle = 1
$ws.Range("A8").Font.Bold = $true

$ws.Range("A9").Value = "CGST"
$ws.Range("A9").Borders.LineStyle = 1
$ws.Range("A9").Font.Bold = $true

$ws.Range("A10").Value = "SGST"
$ws.Range("A10").Borders.LineStyle = 1
$ws.Range("A10").Font.Bold = $true

Write-Host "edit applied"
